$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade record appended as row 4
$ws.Range("A4").Value = 42633.676724537036
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9948
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 19.32
$ws.Range("F4").Value = 19.12
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = -1.04
$ws.Range("I4").Value = $false

# Match the date/time number formatting used by row 3 (column A and G)
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"
